$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

function Add-KeywordRow($stem, $cat) {
    # Seed the new row's formatting from the current last row, then append
    # via the table object so the table range grows and the sheet dimension
    # follows automatically.
    $lastRow = $tbl.ListRows.Item($tbl.ListRows.Count).Range
    $lastRow.Copy()
    $newRow = $tbl.ListRows.Add()
    $newRow.Range.PasteSpecial(-4122)
    $newRow.Range.Cells.Item(1, 1).Value = $stem
    $newRow.Range.Cells.Item(1, 2).Formula = "=LEN(Cluster_Keywords[[#This Row],[Stem]])"
    $newRow.Range.Cells.Item(1, 3).Value = $cat
}

# New keyword stems added to the list.
Add-KeywordRow "Monte" "Childcare"
Add-KeywordRow "Home" "Aged Care"

# "Chef" is renamed to "Chefs".
$bodyRange = $tbl.ListColumns.Item("Stem").DataBodyRange
for ($i = 1; $i -le $bodyRange.Rows.Count; $i++) {
    $cell = $bodyRange.Cells.Item($i, 1)
    if ($cell.Value2 -eq "Chef") {
        $cell.Value = "Chefs"
    }
}

Add-KeywordRow "Meats" "Food"

# Keep the table sorted by Cluster Category then Stem, as it was before.
$tbl.Sort.SortFields.Clear()
$tbl.Sort.SortFields.Add($tbl.ListColumns.Item("Cluster Category").Range)
$tbl.Sort.SortFields.Add($tbl.ListColumns.Item("Stem").Range)
$tbl.Sort.Apply()

# Re-apply the calculated column formula to every row: sorting can leave a
# stray unqualified reference on whichever row lands last, so reassert it
# across the whole column to be safe.
for ($i = 1; $i -le $tbl.ListRows.Count; $i++) {
    $tbl.ListRows.Item($i).Range.Cells.Item(1, 2).Formula = "=LEN(Cluster_Keywords[[#This Row],[Stem]])"
}

[void]$ws.Range("A18").Select()
